$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the booking.com "Locations" test data
$ws.Name = "Locations"

# Populate the location list
$ws.Range("A1").Value = "Locations"
$ws.Range("A2").Value = "New York"
$ws.Range("A3").Value = "California"

# Store the values as text (adds the new numFmtId="49" cell style, applied to
# the populated cells)
$ws.Cells.NumberFormat = "@"

# View tweaks captured in the saved workbook
$excel.ActiveWindow.Zoom = 175
$null = $ws.Range("E8").Select()

# Print setup
$ws.PageSetup.Orientation = 1
